$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.954.78'
$ws.Range("E2").Value = '  +0.32%  '

# Row 3
$ws.Range("D3").Value = '1.818.92'
$ws.Range("E3").Value = '  +0.46%  '

# Row 5
$ws.Range("D5").Value = '''309.93'
$ws.Range("E5").Value = '  +0.10%  '

# Row 6
$ws.Range("E6").Value = '  +0.14%  '

# Row 7
$ws.Range("E7").Value = '  +0.54%  '

# Row 9
$ws.Range("D9").Value = '''0.07356'
$ws.Range("E9").Value = '  -0.13%  '

# Row 10
$ws.Range("D10").Value = '''0.8733'
$ws.Range("E10").Value = '  -0.32%  '

# Row 11
$ws.Range("D11").Value = '''20.28'

# Row 12
$ws.Range("D12").Value = '1.838.86'
$ws.Range("E12").Value = '  +6.05%  '

# Row 13
$ws.Range("D13").Value = '''5.416'
$ws.Range("E13").Value = '  +1.00%  '

# Row 14
$ws.Range("D14").Value = '''0.07114'
$ws.Range("E14").Value = '  +0.90%  '

# Row 15
$ws.Range("D15").Value = '''6.514'
$ws.Range("E15").Value = '  +0.06%  '

# Row 16
$ws.Range("D16").Value = '''91.53'
$ws.Range("E16").Value = '  -0.16%  '

# Row 17
$ws.Range("E17").Value = '  +0.27%  '

# Row 18
$ws.Range("D18").Value = '''0.000008707'

# Row 19
$ws.Range("D19").Value = '''1.002'

# Row 20
$ws.Range("D20").Value = '''14.64'
$ws.Range("E20").Value = '  -0.71%  '

# Row 21
$ws.Range("D21").Value = '26.974.23'
$ws.Range("E21").Value = '  +0.41%  '

# Row 22
$ws.Range("E22").Value = '  -0.24%  '

# Row 23
$ws.Range("E23").Value = '  +0.44%  '

# Row 24
$ws.Range("D24").Value = '2.040.61'
$ws.Range("E24").Value = '  +3.53%  '

# Row 25
$ws.Range("E25").Value = '  -0.29%  '

# Row 26
$ws.Range("D26").Value = '''150.98'
$ws.Range("E26").Value = '  -0.37%  '

# Row 27
$ws.Range("E27").Value = '  +0.06%  '

# Row 28
$ws.Range("D28").Value = '''2.147'
$ws.Range("E28").Value = '  -0.20%  '

# Row 29
$ws.Range("D29").Value = '''5.255'
$ws.Range("E29").Value = '  -1.23%  '

# Row 30
$ws.Range("D30").Value = '''117.41'
$ws.Range("E30").Value = '  +1.39%  '

# Row 31
$ws.Range("D31").Value = '''0.08896'
$ws.Range("E31").Value = '  -0.05%  '

# Row 32
$ws.Range("D32").Value = '''0.7598'
$ws.Range("E32").Value = '  +0.75%  '

# Row 33
$ws.Range("D33").Value = '''1.162'
$ws.Range("E33").Value = '  +0.54%  '

# Row 34
$ws.Range("E34").Value = '  +1.02%  '

# Row 35
$ws.Range("D35").Value = '''2.908'
$ws.Range("E35").Value = '  -0.21%  '

# Row 36
$ws.Range("E36").Value = '  +0.17%  '

# Row 37
$ws.Range("D37").Value = '''1.093'
$ws.Range("E37").Value = '  -0.65%  '

# Row 38
$ws.Range("D38").Value = '''0.05297'
$ws.Range("E38").Value = '  +0.77%  '

# Row 39
$ws.Range("D39").Value = '''0.01948'
$ws.Range("E39").Value = '  -0.83%  '

# Row 40
$ws.Range("D40").Value = '''2.972'
$ws.Range("E40").Value = '  +1.93%  '

# Row 41
$ws.Range("D41").Value = '''7.190'
$ws.Range("E41").Value = '  +0.35%  '

# Row 42
$ws.Range("E42").Value = '  -0.38%  '

# Row 43
$ws.Range("D43").Value = '''2.333'
$ws.Range("E43").Value = '  -4.89%  '

# Row 44
$ws.Range("D44").Value = '''0.1654'
$ws.Range("E44").Value = '  -0.49%  '

# Row 45
$ws.Range("D45").Value = '''8.441'
$ws.Range("E45").Value = '  -0.22%  '

# Row 46
$ws.Range("D46").Value = '''0.4877'
$ws.Range("E46").Value = '  -1.39%  '

# Row 47
$ws.Range("D47").Value = '''10.46'
$ws.Range("E47").Value = '  +1.70%  '

# Row 48
$ws.Range("E48").Value = '  +0.18%  '

# Row 49
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '''1.666'
$ws.Range("E49").Value = '  -0.43%  '

# Row 50
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '''103.37'
$ws.Range("E50").Value = '  +0.13%  '

# Row 51
$ws.Range("D51").Value = '''0.06301'
$ws.Range("E51").Value = '  +0.14%  '

# Clear auto-applied "Text" number-format overrides picked up from the quote-prefix entry,
# so these cells keep the workbook default style (no explicit s= attribute).
$ws.Range("D5").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
